$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The paragraph "Hiển thị các phép đo nhiệt độ, độ ẩm, ánh sáng, pH và CO2)"
# used to be split across two runs with a (now pointless) "_GoBack"
# bookmark sitting between them. Re-typing the same visible text over the
# whole paragraph collapses it back down to a single run and drops the
# bookmark, matching the target XML.
$findText = "Hiển thị các phép đo nhiệt độ, độ ẩm, ánh sáng, pH và CO2)"
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $findText, 2) | Out-Null

# --- Change 2 -----------------------------------------------------------
# Append a new "Prometheus Server" / "Sensor" section (with its bullet
# list) right before the document's trailing empty paragraph, then move
# the "_GoBack" bookmark onto the new final (still empty) paragraph.
$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="4"/><w:bidi w:val="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>Prometheus Server</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:bidi w:val="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="5"/><w:bidi w:val="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Sensor</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">+ List </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr><w:t>all avaiable sensor.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr><w:t>+ Scan all sensor</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr><w:t>+ Collect data from specific sensor (Subcribe/ Unsubcribe)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Times New Roman"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertXML($newBlockXml)
